$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. It belongs right after the
# existing row 13 (chronologically it is inserted before the former row 14),
# so insert a fresh row at position 14 - this pushes the former rows 14..138
# down to 15..139 and grows the sheet's used range to A1:T139.
$ws.Rows(14).Insert()

# Fill in the newly inserted row 14 with the new observation. All the
# "dimension" columns (mercado / producto / categoria / unidad / origen /
# etc.) repeat the constant values used throughout the rest of the sheet.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 45282
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100108
$ws.Range("H14").Value = "Tropicales y subtropicales"
$ws.Range("I14").Value = 100108007
$ws.Range("J14").Value = "Coco"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 30000
$ws.Range("O14").Value = 32000
$ws.Range("P14").Value = 31143
$ws.Range("Q14").Value = "$/malla 20 unidades"
$ws.Range("R14").Value = "Perú"
$ws.Range("S14").Value = 1557
$ws.Range("T14").Value = 20
